$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3920720726908886
$ws.Range("J2").Value = 0.3920720726908886
$ws.Range("M2").Value = 1.275643
$ws.Range("N2").Value = 3.826929
$ws.Range("O2").Value = 0.008652234199457187
$ws.Range("P2").Value = 0.008652234199457187
$ws.Range("Q2").Value = 1.456264268870333
$ws.Range("R2").Value = 13.106378419833
$ws.Range("S2").Value = 0.00339229939598817
$ws.Range("T2").Value = 0.00339229939598817
$ws.Range("I3").Value = 0.3920720726908886
$ws.Range("J3").Value = 0.3920720726908886
$ws.Range("O3").Value = 0.6643867693241158
$ws.Range("P3").Value = 0.6643867693241158
$ws.Range("S3").Value = 0.2604874977173093
$ws.Range("T3").Value = 0.2604874977173093
$ws.Range("I4").Value = 0.3920720726908886
$ws.Range("J4").Value = 0.3920720726908886
$ws.Range("M4").Value = 48.20552666666666
$ws.Range("N4").Value = 144.61658
$ws.Range("O4").Value = 0.326960996476427
$ws.Range("P4").Value = 0.326960996476427
$ws.Range("Q4").Value = 55.03105966696222
$ws.Range("R4").Value = 495.27953700266
$ws.Range("S4").Value = 0.128192275577591
$ws.Range("T4").Value = 0.1281922755775911
$ws.Range("G5").Value = 1.770097666666667
$ws.Range("H5").Value = 5.310293000000001
$ws.Range("I5").Value = 0.6079279273091115
$ws.Range("J5").Value = 0.6079279273091115
$ws.Range("M5").Value = 1.275643
$ws.Range("N5").Value = 3.826929
$ws.Range("O5").Value = 0.008652234199457187
$ws.Range("P5").Value = 0.008652234199457187
$ws.Range("Q5").Value = 2.258012697799667
$ws.Range("R5").Value = 20.322114280197
$ws.Range("S5").Value = 0.005259934803469017
$ws.Range("T5").Value = 0.005259934803469017
$ws.Range("G6").Value = 1.770097666666667
$ws.Range("H6").Value = 5.310293000000001
$ws.Range("I6").Value = 0.6079279273091115
$ws.Range("J6").Value = 0.6079279273091115
$ws.Range("O6").Value = 0.6643867693241158
$ws.Range("P6").Value = 0.6643867693241158
$ws.Range("Q6").Value = 173.3880205736998
$ws.Range("R6").Value = 1560.492185163298
$ws.Range("S6").Value = 0.4038992716068065
$ws.Range("T6").Value = 0.4038992716068065
$ws.Range("G7").Value = 1.770097666666667
$ws.Range("H7").Value = 5.310293000000001
$ws.Range("I7").Value = 0.6079279273091115
$ws.Range("J7").Value = 0.6079279273091115
$ws.Range("M7").Value = 48.20552666666666
$ws.Range("N7").Value = 144.61658
$ws.Range("O7").Value = 0.326960996476427
$ws.Range("P7").Value = 0.326960996476427
$ws.Range("Q7").Value = 85.32849027310445
$ws.Range("R7").Value = 767.9564124579401
$ws.Range("S7").Value = 0.198768720898836
$ws.Range("T7").Value = 0.198768720898836
